# edit.ps1 - apply the "update error in file path for saving converted document"
# commit: refresh the Denver "Maintenance/Plumber" leads sheet with a new scrape
# (new names/companies for existing rows 2-12, plus 14 brand-new rows 13-26).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Step 1 - grow the table downward.
# Column A (s="1": bold font + thin box border + centered) is already
# applied to A2:A15 (14 cells - 11 already populated, 3 still blank).
# Copy that 14-row styled block and paste FORMATS ONLY starting at A13,
# so the 14 new rows (13-26) pick up the existing style index instead of
# minting a brand-new one.
# ------------------------------------------------------------------
$ws.Range("A2:A15").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null   # -4122 = xlPasteFormats
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# Step 2 - refresh rows 2-12 in place: new LinkedIn profile + name,
# new location/title/company per row (index in col A is unchanged).
# ------------------------------------------------------------------

# Row 2
$ws.Range("B2").Formula = '=HYPERLINK("https://www.linkedin.com/in/jack-birdsong-a134883b","jack birdsong")'
$ws.Range("C2").Value = 'Denver'
$ws.Range("D2").Value = 'Journeyman Heating, Ventilation, and Air Conditioning Technician'
$ws.Range("E2").Value = 'Bruce mechanic '

# Row 3
$ws.Range("B3").Formula = '=HYPERLINK("https://www.linkedin.com/in/logan-robinson-608872212","Logan Robinson")'
$ws.Range("C3").Value = 'Denver'
$ws.Range("D3").Value = 'Plumber'
$ws.Range("E3").Value = 'TONY V. PLUMBING & HEATING, INC'

# Row 4
$ws.Range("B4").Formula = '=HYPERLINK("https://www.linkedin.com/in/daniel-a23","Daniel A.")'
$ws.Range("C4").Value = 'Denver'
$ws.Range("D4").Value = 'Heavy Mobile Equipment Mechanic'
$ws.Range("E4").Value = 'United States Space Force'

# Row 5
$ws.Range("B5").Formula = '=HYPERLINK("https://www.linkedin.com/in/greg-curran-918068106","Greg curran")'
$ws.Range("D5").Value = 'Maintenance Technician'
$ws.Range("E5").Value = 'Denver water'

# Row 6
$ws.Range("B6").Formula = '=HYPERLINK("https://www.linkedin.com/in/louie-anglo-jr-b1b492171","Louie Anglo Jr")'
$ws.Range("C6").Value = 'Denver'
$ws.Range("D6").Value = 'Senior plumber'
$ws.Range("E6").Value = 'City and County of Denver'

# Row 7
$ws.Range("B7").Formula = '=HYPERLINK("https://www.linkedin.com/in/brian-mccoy-91054065","Brian McCoy")'
$ws.Range("D7").Value = 'Imaging Service Engineer'
$ws.Range("E7").Value = 'CommonSpirit Health'

# Row 8
$ws.Range("B8").Formula = '=HYPERLINK("https://www.linkedin.com/in/john-eastridge-ii-770864195","John Eastridge II")'
$ws.Range("D8").Value = 'Biomedical Equipment Technician'
$ws.Range("E8").Value = 'RENOVO Solutions'

# Row 9
$ws.Range("B9").Formula = '=HYPERLINK("https://www.linkedin.com/in/julianmeisner","Julian Meisner IV")'
$ws.Range("D9").Value = 'Senior Maintenance Technician'
$ws.Range("E9").Value = 'Camden Property Trust'

# Row 10
$ws.Range("B10").Formula = '=HYPERLINK("https://www.linkedin.com/in/robert-sprague-8a28b0154","Robert Sprague")'
$ws.Range("C10").Value = 'Denver Metropolitan Area'
$ws.Range("D10").Value = 'Plumber'
$ws.Range("E10").Value = 'Roto-Rooter Plumbing and Drain Service'

# Row 11
$ws.Range("B11").Formula = '=HYPERLINK("https://www.linkedin.com/in/bill-howard-043235188","Bill Howard")'
$ws.Range("C11").Value = 'Denver Metropolitan Area'
$ws.Range("D11").Value = 'Residential Plumber'
$ws.Range("E11").Value = 'Canyon Plumbing'

# Row 12
$ws.Range("B12").Formula = '=HYPERLINK("https://www.linkedin.com/in/kevin-bertram-render-9731b6151","Kevin Bertram-Render")'
$ws.Range("C12").Value = 'Denver'
$ws.Range("D12").Value = 'Apprentice Plumber'
$ws.Range("E12").Value = 'Bates Mechanical, Inc.'

# ------------------------------------------------------------------
# Step 3 - append the 14 brand-new rows 13-26 (index 11-24).
# ------------------------------------------------------------------

# Row 13
$ws.Range("A13").Value = 11
$ws.Range("B13").Formula = '=HYPERLINK("https://www.linkedin.com/in/chris-sing-013946142","Chris Sing")'
$ws.Range("C13").Value = 'Denver'
$ws.Range("D13").Value = 'Master Plumber'
$ws.Range("E13").Value = 'Applewood Plumbing Heating & Electric'

# Row 14
$ws.Range("A14").Value = 12
$ws.Range("B14").Formula = '=HYPERLINK("https://www.linkedin.com/in/andrew-romero-70052415b","Andrew Romero")'
$ws.Range("C14").Value = 'Englewood'
$ws.Range("D14").Value = 'Maintenance Technician'
$ws.Range("E14").Value = 'Meadow Gold Dairies'

# Row 15
$ws.Range("A15").Value = 13
$ws.Range("B15").Formula = '=HYPERLINK("https://www.linkedin.com/in/jasonbandykarma","Jason Bandy")'
$ws.Range("C15").Value = 'Denver'
$ws.Range("D15").Value = 'Plumber'
$ws.Range("E15").Value = 'Searching'

# Row 16
$ws.Range("A16").Value = 14
$ws.Range("B16").Formula = '=HYPERLINK("https://www.linkedin.com/in/lue-lor-97a725159","Lue Lor")'
$ws.Range("C16").Value = 'Denver'
$ws.Range("D16").Value = 'Maintenance Technician'
$ws.Range("E16").Value = 'Medtronic'

# Row 17
$ws.Range("A17").Value = 15
$ws.Range("B17").Formula = '=HYPERLINK("https://www.linkedin.com/in/darvi-olivares-097b39135","Darvi Olivares")'
$ws.Range("C17").Value = 'Littleton'
$ws.Range("D17").Value = 'Maintenance Technician'
$ws.Range("E17").Value = 'MHCD'

# Row 18
$ws.Range("A18").Value = 16
$ws.Range("B18").Formula = '=HYPERLINK("https://www.linkedin.com/in/alex-montgomery-67448a15b","Alex Montgomery")'
$ws.Range("C18").Value = 'Littleton'
$ws.Range("D18").Value = 'Sales Support Representative'
$ws.Range("E18").Value = 'AEE Solar'

# Row 19
$ws.Range("A19").Value = 17
$ws.Range("B19").Formula = '=HYPERLINK("https://www.linkedin.com/in/brandon-sweet-1aa262a9","Brandon Sweet")'
$ws.Range("C19").Value = 'Littleton'
$ws.Range("D19").Value = 'Plumber'
$ws.Range("E19").Value = 'Wheatridge Plumbing & Heating'

# Row 20
$ws.Range("A20").Value = 18
$ws.Range("B20").Formula = '=HYPERLINK("https://www.linkedin.com/in/josue-flores-3250051a0","Josue Flores")'
$ws.Range("C20").Value = 'Denver'
$ws.Range("D20").Value = 'Plumber'
$ws.Range("E20").Value = 'Trautman & Shreve, Inc'

# Row 21
$ws.Range("A21").Value = 19
$ws.Range("B21").Formula = '=HYPERLINK("https://www.linkedin.com/in/connor-brady-575a9189","Connor Brady")'
$ws.Range("C21").Value = 'Denver'
$ws.Range("D21").Value = 'Plumber'
$ws.Range("E21").Value = 'AAA Service Plumbing, Heating, and Electric'

# Row 22
$ws.Range("A22").Value = 20
$ws.Range("B22").Formula = '=HYPERLINK("https://www.linkedin.com/in/joseph-calimpong-70574969","Joseph Calimpong")'
$ws.Range("C22").Value = 'Denver Metropolitan Area'
$ws.Range("D22").Value = 'Carwash Maintenance Tech'
$ws.Range("E22").Value = 'Hi Performance Car Wash'

# Row 23
$ws.Range("A23").Value = 21
$ws.Range("B23").Formula = '=HYPERLINK("https://www.linkedin.com/in/tyler-hudziec-b322271a2","Tyler Hudziec")'
$ws.Range("C23").Value = 'Denver Metropolitan Area'
$ws.Range("D23").Value = 'Heating Air Conditioning Specialist'
$ws.Range("E23").Value = 'Denver Public Library'

# Row 24
$ws.Range("A24").Value = 22
$ws.Range("B24").Formula = '=HYPERLINK("https://www.linkedin.com/in/deon-g-575bb7245","Deon G.")'
$ws.Range("C24").Value = 'Denver'

# Row 25
$ws.Range("A25").Value = 23
$ws.Range("B25").Formula = '=HYPERLINK("https://www.linkedin.com/in/tommy-carrillo-258b00245","Tommy Carrillo")'
$ws.Range("C25").Value = 'Denver'
$ws.Range("D25").Value = 'Maintenance Technician'
$ws.Range("E25").Value = 'Property Maintenance & Management Services Ltd.'

# Row 26
$ws.Range("A26").Value = 24
$ws.Range("B26").Formula = '=HYPERLINK("https://www.linkedin.com/in/jose-saenz-131a96193","Jose Saenz")'
$ws.Range("C26").Value = 'Denver Metropolitan Area'
$ws.Range("D26").Value = 'Licensed apprentice plumber'
$ws.Range("E26").Value = 'CHS Plumbing'
